# Update cryptos list Price (D) and Volume(1h) (E) values per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.050.29"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.299.87"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D5").Value = "'300.12"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'97.95"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  +2.42%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.515"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'36.12"
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "'17.74"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").Value = "'6.87"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "2.661.11"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "2.297.71"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "'0.787"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "42.957.94"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'12.79"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "'6.14"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'68.73"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").Value = "'237.83"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "'24.97"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "'164.59"
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "'9.13"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'33.01"
$ws.Range("E32").Value = "  -4.44%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("D35").Value = "'4.80"
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("D36").Value = "'17.99"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "'0.0697"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").Value = "'2.78"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "2.017.00"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").Value = "'0.0286"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").Value = "'2.20"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("D46").Value = "'10.41"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("D47").Value = "'17.43"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "'2.83"
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("D49").Value = "'54.10"
$ws.Range("E49").Value = "  -3.06%  "
$ws.Range("D50").Value = "2.529.73"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").Value = "'1.53"
$ws.Range("E51").Value = "  -1.57%  "
